$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 238.71428
$ws.Range("I5").Value = 375
$ws.Range("J5").Value = 57
$ws.Range("K5").Value = 375
$ws.Range("L5").Value = 57
$ws.Range("M5").Value = -260
$ws.Range("N5").Value = -287

# Row 33
$ws.Range("H33").Value = 131.33333
$ws.Range("I33").Value = 131.33333
$ws.Range("K33").Value = 131.33333
$ws.Range("M33").Value = 97.66667000000001

# Row 100
$ws.Range("H100").Value = 1939.2
$ws.Range("I100").Value = 1924.25
$ws.Range("J100").Value = 1999
$ws.Range("K100").Value = 1924.25
$ws.Range("L100").Value = 1999
$ws.Range("M100").Value = -1383.25
$ws.Range("N100").Value = -3081

# Row 113
$ws.Range("H113").Value = 33320.715
$ws.Range("I113").Value = 24849.75
$ws.Range("J113").Value = 44615.332
$ws.Range("K113").Value = 24849.75
$ws.Range("L113").Value = 44615.332
$ws.Range("M113").Value = -21595.75
$ws.Range("N113").Value = -51123.332

# Row 137
$ws.Range("H137").Value = 2286.3333
$ws.Range("I137").Value = 2149.6667
$ws.Range("K137").Value = 6449.000100000001
$ws.Range("M137").Value = -3899.000100000001

# Row 141
$ws.Range("H141").Value = 3730
$ws.Range("I141").Value = 2595
$ws.Range("K141").Value = 7785
$ws.Range("M141").Value = -2605

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 13
$ws.Range("I2").Value = 13
$ws.Range("K2").Value = 13
$ws.Range("M2").Value = 100

# Row 61
$ws.Range("H61").Value = 1704.4
$ws.Range("I61").Value = 1704.4
$ws.Range("K61").Value = 1704.4
$ws.Range("M61").Value = -1492.4

# Row 116
$ws.Range("H116").Value = 13
$ws.Range("I116").Value = 13
$ws.Range("K116").Value = 13
$ws.Range("M116").Value = 2281

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# Row 136
$ws.Range("H136").Value = 1704.4
$ws.Range("I136").Value = 1704.4
$ws.Range("K136").Value = 5113.200000000001
$ws.Range("M136").Value = -2563.200000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 13
$ws.Range("K3").Value = 13
$ws.Range("M3").Value = 101

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8545.1
$ws.Range("I31").Value = 8316.333000000001
$ws.Range("K31").Value = 8316.333000000001
$ws.Range("M31").Value = -8021.333000000001

# Row 34
$ws.Range("H34").Value = 8545.1
$ws.Range("I34").Value = 8316.333000000001
$ws.Range("K34").Value = 8316.333000000001
$ws.Range("M34").Value = -8114.333000000001

# Row 39
$ws.Range("H39").Value = 11810.4
$ws.Range("J39").Value = 23664.25
$ws.Range("L39").Value = 23664.25
$ws.Range("N39").Value = -24446.25

# Row 49
$ws.Range("H49").Value = 11810.4
$ws.Range("J49").Value = 23664.25
$ws.Range("L49").Value = 23664.25
$ws.Range("N49").Value = -24028.25

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

# Row 99
$ws.Range("H99").Value = 836268.5
$ws.Range("I99").Value = 557691.3
$ws.Range("K99").Value = 557691.3
$ws.Range("M99").Value = -556193.3

# Row 107
$ws.Range("H107").Value = 1006.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1006.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1006.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4846.5

# Row 126
$ws.Range("H126").Value = 836268.5
$ws.Range("I126").Value = 557691.3
$ws.Range("K126").Value = 1673073.9
$ws.Range("M126").Value = -1670603.9

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 3000
$ws.Range("J39").Value = 3000
$ws.Range("L39").Value = 9000
$ws.Range("N39").Value = -9588

# Row 55
$ws.Range("H55").Value = 1915.4166
$ws.Range("J55").Value = 5066.6665
$ws.Range("L55").Value = 15199.9995
$ws.Range("N55").Value = -15553.9995

# Row 97
$ws.Range("H97").Value = 498.5
$ws.Range("I97").Value = 498.5
$ws.Range("K97").Value = 1495.5
$ws.Range("M97").Value = -999.5

# Row 129
$ws.Range("H129").Value = 2163.3333
$ws.Range("I129").Value = 1990
$ws.Range("J129").Value = 2250
$ws.Range("K129").Value = 5970
$ws.Range("L129").Value = 6750
$ws.Range("M129").Value = -970
$ws.Range("N129").Value = -16750

# Row 132
$ws.Range("H132").Value = 10200
$ws.Range("I132").Value = 8333.333000000001
$ws.Range("J132").Value = 15800
$ws.Range("K132").Value = 74999.997
$ws.Range("L132").Value = 142200
$ws.Range("M132").Value = -72469.997
$ws.Range("N132").Value = -147260

# Row 140
$ws.Range("H140").Value = 429
$ws.Range("I140").Value = 429
$ws.Range("K140").Value = 1287
$ws.Range("M140").Value = 3893

# Row 141
$ws.Range("H141").Value = 1843.3334
$ws.Range("I141").Value = 1843.3334
$ws.Range("K141").Value = 5530.0002
$ws.Range("M141").Value = -350.0002000000004

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1532.3334
$ws.Range("I97").Value = 1439
$ws.Range("J97").Value = 1999
$ws.Range("K97").Value = 1439
$ws.Range("L97").Value = 1999
$ws.Range("M97").Value = -943
$ws.Range("N97").Value = -2991

# Row 102
$ws.Range("H102").Value = 399999
$ws.Range("I102").Value = 399999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 399999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -398377
$ws.Range("N102").ClearContents()

# Row 126
$ws.Range("H126").Value = 8217.182000000001
$ws.Range("I126").Value = 8271.286
$ws.Range("K126").Value = 24813.858
$ws.Range("M126").Value = -22343.858

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 828.5714
$ws.Range("I16").Value = 828.5714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 828.5714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -658.5714
$ws.Range("N16").ClearContents()

# Row 46
$ws.Range("H46").Value = 501499
$ws.Range("I46").Value = 1000750
$ws.Range("J46").Value = 2248
$ws.Range("K46").Value = 1000750
$ws.Range("L46").Value = 2248
$ws.Range("M46").Value = -1000562
$ws.Range("N46").Value = -2624

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 600
$ws.Range("I34").Value = 600
$ws.Range("K34").Value = 600
$ws.Range("M34").Value = -397

# Row 54
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 30000
$ws.Range("K54").Value = 30000
$ws.Range("M54").Value = -29480

# Row 126
$ws.Range("H126").Value = 1254.421
$ws.Range("I126").Value = 1229.6666
$ws.Range("K126").Value = 3688.9998
$ws.Range("M126").Value = -1218.9998

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
